$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted at row 17 (pushing the existing
# rows 17-52 down to 18-53). Insert a blank row at 18, duplicate the
# (still unmodified) row 17 into it, then overwrite row 17 with the new
# record's values.
$ws.Rows(18).Insert()
$ws.Range("A17:T17").Copy($ws.Range("A18:T18"))

$ws.Range("D17").Value = 44459
$ws.Range("M17").Value = 50
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 20000
$ws.Range("S17").Value = 2000
